$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.381.46"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "2.620.59"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'308.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("D6").Value = "'99.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.20%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +3.15%  "
$ws.Range("D10").Value = "'39.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.43%  "
$ws.Range("E11").Value = "  +1.99%  "
$ws.Range("D12").Value = "'54.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").Value = "'8.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.31%  "
$ws.Range("D14").Value = "3.021.46"
$ws.Range("E14").Value = "  +1.27%  "
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("D16").Value = "2.611.61"
$ws.Range("E16").Value = "  +0.92%  "
$ws.Range("E17").Value = "  +3.28%  "
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "46.568.97"
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("E20").Value = "  +1.78%  "
$ws.Range("D21").Value = "'13.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.98%  "
$ws.Range("E22").Value = "  +3.91%  "
$ws.Range("E23").Value = "  +2.98%  "
$ws.Range("D24").Value = "'276.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.87%  "
$ws.Range("E25").Value = "  +2.31%  "
$ws.Range("D26").Value = "'2.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.27%  "
$ws.Range("D27").Value = "'29.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +12.63%  "
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("E29").Value = "  -1.20%  "
$ws.Range("D30").Value = "'10.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.82%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "'38.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.39%  "
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "'2.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.54%  "
$ws.Range("D33").Value = "'6.47"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +9.11%  "
$ws.Range("E34").Value = "  -3.80%  "
$ws.Range("E35").Value = "  +2.09%  "
$ws.Range("E36").Value = "  -3.16%  "
$ws.Range("E37").Value = "  +0.93%  "
$ws.Range("D38").Value = "'152.29"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.45%  "
$ws.Range("E39").Value = "  +2.52%  "
$ws.Range("D40").Value = "'0.123"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.49%  "
$ws.Range("D41").Value = "'24.34"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +38.65%  "
$ws.Range("D42").Value = "'16.02"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.54%  "
$ws.Range("D43").Value = "'0.0330"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.44%  "
$ws.Range("E44").Value = "  +0.82%  "
$ws.Range("D45").Value = "'4.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.53%  "
$ws.Range("D46").Value = "2.142.11"
$ws.Range("E46").Value = "  +5.22%  "
$ws.Range("D47").Value = "'0.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").Value = "'95.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.86%  "
$ws.Range("E49").Value = "  +8.17%  "
$ws.Range("D50").Value = "'109.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.66%  "
$ws.Range("E51").Value = "  -3.29%  "
